# "agr fechou, com tooltip e tudo" -- add the AFN->AFD conversion result
# table (Sheet3), wire up the "CONVERTIDO" labels/tooltip button on
# Sheet1/Sheet2, and select the new sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet1: a new row is inserted above the old "AFD" example block,
# holding a "CONVERTIDO" label.
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws1.Rows.Item(15).Insert()
$ws1.Cells.Item(15, 1).Value = "CONVERTIDO"
$ws1.Range("F9").Select()

# ---------------------------------------------------------------
# Sheet2: widen column E, label it "CONVERTIDO >>", and clear the
# big ">>" button's leftover text (keeps its style).
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Columns.Item(5).ColumnWidth = 22.59
$ws2.Cells.Item(5, 5).Value = "CONVERTIDO >>"
$ws2.Cells.Item(6, 6).Value = ""

# ---------------------------------------------------------------
# Sheet3 (new): the converted AFD table + the transition table.
# ---------------------------------------------------------------
$ws3 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wb.Worksheets.Item($wb.Worksheets.Count))
$ws3.Name = "Sheet3"

$ws3.Cells.Item(1, 2).Value = "a"
$ws3.Cells.Item(1, 3).Value = "b"
$ws3.Cells.Item(1, 4).Value = "c"
$ws3.Cells.Item(1, 5).Value = "@"
$ws3.Cells.Item(1, 10).Value = "a"
$ws3.Cells.Item(1, 11).Value = "b"
$ws3.Cells.Item(1, 12).Value = "c"

$ws3.Cells.Item(2, 1).Value = "q0"
$ws3.Cells.Item(2, 2).Value = "q0"
$ws3.Cells.Item(2, 3).Value = "q0"
$ws3.Cells.Item(2, 4).Value = "q0"
$ws3.Cells.Item(2, 5).Value = "q1,q2,q4"
$ws3.Cells.Item(2, 9).Value = "q0q1q2q4"
$ws3.Cells.Item(2, 10).Value = "q0q1q2q4qf"
$ws3.Cells.Item(2, 11).Value = "q0q1q2q4q3"
$ws3.Cells.Item(2, 12).Value = "q0q1q2q4q5"

$ws3.Cells.Item(3, 1).Value = "q1"
$ws3.Cells.Item(3, 2).Value = "qf"
$ws3.Cells.Item(3, 9).Value = "q0q1q2q4qf"
$ws3.Cells.Item(3, 10).Value = "q0q1q2q4"
$ws3.Cells.Item(3, 11).Value = "q0q1q2q4"
$ws3.Cells.Item(3, 12).Value = "q0q1q2q4"

$ws3.Cells.Item(4, 1).Value = "q2"
$ws3.Cells.Item(4, 3).Value = "q3"
$ws3.Cells.Item(4, 9).Value = "q0q1q2q4q3"
$ws3.Cells.Item(4, 10).Value = "q0q1q2q4qf"
$ws3.Cells.Item(4, 11).Value = "q0q1q2q4q3qf"
$ws3.Cells.Item(4, 12).Value = "q0q1q2q4q5"

$ws3.Cells.Item(5, 1).Value = "q3"
$ws3.Cells.Item(5, 3).Value = "qf"
$ws3.Cells.Item(5, 7).Value = "CONVERTIDO >>"
$ws3.Cells.Item(5, 9).Value = "q0q1q2q4q5"
$ws3.Cells.Item(5, 10).Value = "q0q1q2q4qf"
$ws3.Cells.Item(5, 11).Value = "q0q1q2q4q3"
$ws3.Cells.Item(5, 12).Value = "q0q1q2q4q5q6"

$ws3.Cells.Item(6, 1).Value = "q4"
$ws3.Cells.Item(6, 4).Value = "q5"
$ws3.Cells.Item(6, 9).Value = "q0q1q2q4q3qf"
$ws3.Cells.Item(6, 10).Value = "q0q1q2q4qf"
$ws3.Cells.Item(6, 11).Value = "q0q1q2q4q3qf"
$ws3.Cells.Item(6, 12).Value = "q0q1q2q4q5"

$ws3.Cells.Item(7, 1).Value = "q5"
$ws3.Cells.Item(7, 4).Value = "q6"
$ws3.Cells.Item(7, 9).Value = "q0q1q2q4q5q6"
$ws3.Cells.Item(7, 10).Value = "q0q1q2q4qf"
$ws3.Cells.Item(7, 11).Value = "q0q1q2q4q3"
$ws3.Cells.Item(7, 12).Value = "q0q1q2q4q5q6qf"

$ws3.Cells.Item(8, 1).Value = "q6"
$ws3.Cells.Item(8, 4).Value = "qf"
$ws3.Cells.Item(8, 9).Value = "q0q1q2q4q5q6qf"
$ws3.Cells.Item(8, 10).Value = "q0q1q2q4qf"
$ws3.Cells.Item(8, 11).Value = "q0q1q2q4q3"
$ws3.Cells.Item(8, 12).Value = "q0q1q2q4q5q6qf"

$ws3.Cells.Item(9, 1).Value = "qf"

$ws3.Columns.Item(9).ColumnWidth = 14.92
$ws3.Columns.Item(10).ColumnWidth = 10.59
$ws3.Columns.Item(11).ColumnWidth = 12.76
$ws3.Columns.Item(12).ColumnWidth = 13.09

$ws3.Range("I8").Select()
